# Auto-generated edit script: updates Leve profit-calculation columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 325.46155
$ws.Range("I11").Value = 325.46155
$ws.Range("K11").Value = 325.46155
$ws.Range("M11").Value = -185.46155
$ws.Range("H17").Value = 3772.4482
$ws.Range("J17").Value = 3896
$ws.Range("L17").Value = 11688
$ws.Range("N17").Value = -12024
$ws.Range("H127").Value = 2030.4166
$ws.Range("I127").Value = 1104.8
$ws.Range("K127").Value = 3314.4
$ws.Range("M127").Value = 1645.6
$ws.Range("H132").Value = 6551.25
$ws.Range("I132").Value = 7068.3335
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 21205.0005
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -18675.0005
$ws.Range("N132").Value = -20060
$ws.Range("H138").Value = 1209.75
$ws.Range("I138").Value = 419.5
$ws.Range("K138").Value = 1258.5
$ws.Range("M138").Value = 3881.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12058.077
$ws.Range("I32").Value = 10562.917
$ws.Range("K32").Value = 10562.917
$ws.Range("M32").Value = -10275.917
$ws.Range("H37").Value = 25000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H55").Value = 23694.46
$ws.Range("J55").Value = 26816.363
$ws.Range("L55").Value = 26816.363
$ws.Range("N55").Value = -27446.363
$ws.Range("H80").Value = 39997.777
$ws.Range("J80").Value = 39997.777
$ws.Range("L80").Value = 39997.777
$ws.Range("N80").Value = -41993.777
$ws.Range("H83").Value = 39997.777
$ws.Range("J83").Value = 39997.777
$ws.Range("L83").Value = 119993.331
$ws.Range("N83").Value = -129977.331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 241
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H86").Value = 2079.1428
$ws.Range("I86").Value = 2079.1428
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2079.1428
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -956.1428000000001
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2079.1428
$ws.Range("I89").Value = 2079.1428
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10395.714
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4779.714
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 29997.5
$ws.Range("J50").Value = 29997.5
$ws.Range("L50").Value = 29997.5
$ws.Range("N50").Value = -31247.5
$ws.Range("H60").Value = 18938.3
$ws.Range("I60").Value = 9848.25
$ws.Range("K60").Value = 9848.25
$ws.Range("M60").Value = -9337.25
$ws.Range("H74").Value = 37920.54
$ws.Range("J74").Value = 38039.418
$ws.Range("L74").Value = 38039.418
$ws.Range("N74").Value = -39787.418
$ws.Range("H77").Value = 37920.54
$ws.Range("J77").Value = 38039.418
$ws.Range("L77").Value = 114118.254
$ws.Range("N77").Value = -122854.254
$ws.Range("H105").Value = 1450
$ws.Range("I105").Value = 1450
$ws.Range("K105").Value = 1450
$ws.Range("M105").Value = 297

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 883.36365
$ws.Range("I11").Value = 792.8333
$ws.Range("J11").Value = 992
$ws.Range("K11").Value = 2378.4999
$ws.Range("L11").Value = 2976
$ws.Range("M11").Value = -2238.4999
$ws.Range("N11").Value = -3256
$ws.Range("H23").Value = 505.7143
$ws.Range("I23").Value = 625
$ws.Range("J23").Value = 416.25
$ws.Range("K23").Value = 1875
$ws.Range("L23").Value = 1248.75
$ws.Range("M23").Value = -1640
$ws.Range("N23").Value = -1718.75
$ws.Range("H34").Value = 466.66666
$ws.Range("I34").Value = 450
$ws.Range("J34").Value = 500
$ws.Range("K34").Value = 1350
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -1266
$ws.Range("N34").Value = -1668
$ws.Range("H39").Value = 9333.333000000001
$ws.Range("I39").Value = 8000
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 24000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -23706
$ws.Range("N39").Value = -30588
$ws.Range("H55").Value = 1937.25
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 75000
$ws.Range("I74").Value = 75000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 225000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -223939
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 75000
$ws.Range("I77").Value = 75000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 675000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -669696
$ws.Range("N77").ClearContents()
$ws.Range("H113").Value = 677.8
$ws.Range("H122").Value = 626
$ws.Range("I122").Value = 598.4
$ws.Range("K122").Value = 5385.599999999999
$ws.Range("M122").Value = -2935.599999999999
$ws.Range("H129").Value = 78.8
$ws.Range("I129").Value = 74.75
$ws.Range("J129").Value = 95
$ws.Range("K129").Value = 224.25
$ws.Range("L129").Value = 285
$ws.Range("M129").Value = 4775.75
$ws.Range("N129").Value = -10285
$ws.Range("H137").Value = 1000
$ws.Range("J137").Value = 1000
$ws.Range("L137").Value = 3000
$ws.Range("N137").Value = -13200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 12900
$ws.Range("H57").Value = 25250
$ws.Range("J57").Value = 25250
$ws.Range("L57").Value = 25250
$ws.Range("N57").Value = -26890
$ws.Range("H80").Value = 3465
$ws.Range("J80").Value = 3998
$ws.Range("L80").Value = 3998
$ws.Range("N80").Value = -5994
$ws.Range("H83").Value = 3465
$ws.Range("J83").Value = 3998
$ws.Range("L83").Value = 19990
$ws.Range("N83").Value = -29974

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2357.1428
$ws.Range("J22").Value = 3250
$ws.Range("L22").Value = 3250
$ws.Range("N22").Value = -3840
$ws.Range("H27").Value = 2357.1428
$ws.Range("J27").Value = 3250
$ws.Range("L27").Value = 3250
$ws.Range("N27").Value = -3464
$ws.Range("H46").Value = 4666.6665
$ws.Range("J46").Value = 4666.6665
$ws.Range("L46").Value = 4666.6665
$ws.Range("N46").Value = -5042.6665
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 374
$ws.Range("J9").Value = 374
$ws.Range("L9").Value = 374
$ws.Range("N9").Value = -654
$ws.Range("H136").Value = 1604.1177
$ws.Range("I136").Value = 1573.7858
$ws.Range("K136").Value = 4721.357400000001
$ws.Range("M136").Value = -2171.357400000001
